# Weekly update: insert a new price record for "Terminal La Palmera de La
# Serena - Zanahoria" at row 472 (new date 45077), pushing the existing
# rows 472:538 down to 473:539.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 472; this shifts rows 472:538
# down to 473:539 and expands the used range to A1:R539.
$ws.Rows("472:472").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A472").Value2 = 8
$ws.Range("B472").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C472").Value2 = "Coquimbo"
$ws.Range("D472").Value2 = 45077
$ws.Range("E472").Value2 = 4
$ws.Range("F472").Value2 = 100114013
$ws.Range("G472").Value2 = "Zanahoria"
$ws.Range("H472").Value2 = "Sin especificar"
$ws.Range("I472").Value2 = "Primera"
$ws.Range("J472").Value2 = 480
$ws.Range("K472").Value2 = 5000
$ws.Range("L472").Value2 = 6000
$ws.Range("M472").Value2 = 5500
$ws.Range("N472").Value2 = "$/saco 20 kilos"
$ws.Range("O472").Value2 = "Provincia del Elquí"
$ws.Range("P472").Value2 = 275
$ws.Range("Q472").Value2 = 20
$ws.Range("R472").Value2 = "Hortaliza"
